$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDescription = "O Plano de Gestão Proativa de Seca do Hidrossistema Patu, concluído em 2022, foi o primeiro do Ceará e do Brasil, desenvolvido no âmbito do Programa Cientista Chefe – Recursos Hídricos. Criado de forma coletiva e interdisciplinar, envolveu universidades, colegiados de participação social (Comitê de Bacia do Banabuiú e Comissão Gestora do Açude Patu) e gestores públicos, servindo de modelo para outros planos no estado.
Localizado no município de Senador Pompeu, no Sertão Central do Ceará, o hidrossistema integra a Região Hidrográfica do Banabuiú, tendo como principal reservatório o Açude Patu, construído em 1987 pelo DNOCS, com capacidade de 65,1 milhões de m³ e área de drenagem de 993,5 km². É usado para abastecimento humano, irrigação, dessedentação animal e agropecuária. A seca de 2012 reduziu drasticamente sua disponibilidade hídrica, afetando também Milhã e Quixeramobim.
O Sistema de Suporte à Decisão do plano oferece uma interface simples para monitorar e divulgar os estados de seca, as ações previstas para cada cenário e registrar novas medidas, atendendo especialmente à Câmara Técnica Permanente de Gestão Proativa de Seca da Região e à Gerência Regional da Cogerh do Banabuiú."

$ws.Range("A2").Value = $newDescription

$ws.Rows.Item(2).RowHeight = 173.25

$ws.Range("A2").Select()
